# The resume's table "header" row (S.No / ... ) was originally placed as the
# LAST row of each table; it needs to become the FIRST row instead.
# We recreate it as a brand-new first row (so formatting/order matches a
# real "move to top"), fill in its text, then delete the old trailing
# header row.

$d = $word.ActiveDocument

function Move-HeaderRowToTop($table, [string[]]$headers) {
    $newRow = $table.Rows.Add($table.Rows(1))
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $newRow.Cells($i + 1).Range.Text = $headers[$i]
    }

    # The original trailing header row shifted down by one position (to the
    # new last index) once the fresh row was inserted at the top; fetch it
    # fresh rather than relying on a pre-insert reference, since that
    # becomes stale after the row collection shifts.
    $table.Rows($table.Rows.Count).Delete()
}

Move-HeaderRowToTop $d.Tables(1) @("S.No", "Institution", "Degree", "Year")
Move-HeaderRowToTop $d.Tables(2) @("S.No", "Company & Role", "Time Period", "Description")
Move-HeaderRowToTop $d.Tables(3) @("S.No", "Project Name", "Technologies", "Description")
